$p = $ppt.ActivePresentation

# The commit swaps the OOXML content of ppt/theme/theme1.xml (the theme
# wired to the slide master, originally the "Integral" theme) with
# ppt/theme/theme2.xml (the theme wired to the notes master, originally
# the default "Office Theme"). fontScheme/fmtScheme are identical between
# the two themes already, so the only observable difference is the
# <a:clrScheme> color values (and names). The only real, writable surface
# for theme colors exposed by this host is Slide.ThemeColorScheme, which
# maps 1:1 onto ppt/theme/theme1.xml's <a:clrScheme> slots in document
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      (hex 44546A, stored as BGR int)
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      (hex E7E6E6)
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  (hex 5B9BD5)
$tcs.Colors(6).RGB  = 0x317DED   # accent2  (hex ED7D31)
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  (hex A5A5A5)
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  (hex FFC000)
$tcs.Colors(9).RGB  = 0xC47244   # accent5  (hex 4472C4)
$tcs.Colors(10).RGB = 0x47AD70   # accent6  (hex 70AD47)
$tcs.Colors(11).RGB = 0xC16305   # hlink    (hex 0563C1)
$tcs.Colors(12).RGB = 0x724F95   # folHlink (hex 954F72)
